$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.005699780259004132
$ws.Range("C2").Value = 0.004815947175593493
$ws.Range("D2").Value = 0.006019207552184639
$ws.Range("E2").Value = 0.00695184952048763
$ws.Range("B3").Value = 1.739704219956449
$ws.Range("C3").Value = 2.089661935040025
$ws.Range("D3").Value = 3.093403766648005
$ws.Range("E3").Value = 3.804266060556011
$ws.Range("B4").Value = -0.003397341679187313
$ws.Range("C4").Value = -0.005060226800726117
$ws.Range("D4").Value = -0.005880069239839422
$ws.Range("E4").Value = -0.005935788986519947
$ws.Range("B5").Value = -1.533100060524297
$ws.Range("C5").Value = -2.730250707311016
$ws.Range("D5").Value = -3.838105272924532
$ws.Range("E5").Value = -4.513695152992129
$ws.Range("B6").Value = -0.0001868963237655583
$ws.Range("C6").Value = 0.0009864865843146812
$ws.Range("D6").Value = -0.00559623650165808
$ws.Range("E6").Value = -0.001881410720966518
$ws.Range("B7").Value = -0.16805332308191
$ws.Range("C7").Value = 0.4010659441796154
$ws.Range("D7").Value = -1.807871582119023
$ws.Range("E7").Value = -0.5492005480496613
$ws.Range("B8").Value = 0.004785687961972761
$ws.Range("C8").Value = 0.00452641617384427
$ws.Range("D8").Value = 0.005986258223583007
$ws.Range("E8").Value = 0.007600139250810986
$ws.Range("B9").Value = 1.457869192700978
$ws.Range("C9").Value = 1.95993084498245
$ws.Range("D9").Value = 2.976728176949976
$ws.Range("E9").Value = 3.897709540138819
$ws.Range("B10").Value = -0.00816532825476327
$ws.Range("C10").Value = -0.007032806108901811
$ws.Range("D10").Value = -0.007338361278754536
$ws.Range("E10").Value = -0.005928761921112338
$ws.Range("B11").Value = -2.885206247699012
$ws.Range("C11").Value = -3.458849286189286
$ws.Range("D11").Value = -4.756818425555082
$ws.Range("E11").Value = -4.871914683516339
$ws.Range("B12").Value = -0.002162129197188976
$ws.Range("C12").Value = -0.005996239352825823
$ws.Range("D12").Value = -0.006359798070396438
$ws.Range("E12").Value = -0.001448461190821408
$ws.Range("B13").Value = -1.538323721476352
$ws.Range("C13").Value = -2.061491092116992
$ws.Range("D13").Value = -2.163517623739851
$ws.Range("E13").Value = -0.4309032168234385
$ws.Range("B14").Value = 0.005055430241886743
$ws.Range("C14").Value = 0.004683093404372679
$ws.Range("D14").Value = 0.006187681801526428
$ws.Range("E14").Value = 0.007856985027583507
$ws.Range("B15").Value = 1.521471578969507
$ws.Range("C15").Value = 1.978462176213194
$ws.Range("D15").Value = 3.018251838362525
$ws.Range("E15").Value = 3.969858075212485
$ws.Range("B16").Value = -0.006401444703493888
$ws.Range("C16").Value = -0.007379937413733258
$ws.Range("D16").Value = -0.007351262216361159
$ws.Range("E16").Value = -0.005638523698822219
$ws.Range("B17").Value = -2.071757768678078
$ws.Range("C17").Value = -3.606705359259601
$ws.Range("D17").Value = -4.840882952723406
$ws.Range("E17").Value = -4.786377614314292
$ws.Range("B18").Value = -0.001762858819559474
$ws.Range("C18").Value = -0.005698897265180042
$ws.Range("D18").Value = -0.005823424256875306
$ws.Range("E18").Value = -0.00235109639791297
$ws.Range("B19").Value = -1.253030264629352
$ws.Range("C19").Value = -2.042743421186345
$ws.Range("D19").Value = -2.02158949058392
$ws.Range("E19").Value = -0.7436538190882448
$ws.Range("B20").Value = 0.005987230815273466
$ws.Range("C20").Value = 0.005528347939844945
$ws.Range("D20").Value = 0.007016831334659212
$ws.Range("E20").Value = 0.008495335097801295
$ws.Range("B21").Value = 1.787900845620036
$ws.Range("C21").Value = 2.291499481653046
$ws.Range("D21").Value = 3.387198153959176
$ws.Range("E21").Value = 4.388967040534051
$ws.Range("B22").Value = -0.00659235594305914
$ws.Range("C22").Value = -0.008101624932892503
$ws.Range("D22").Value = -0.007750688364177777
$ws.Range("E22").Value = -0.006350500135489026
$ws.Range("B23").Value = -2.138047777371845
$ws.Range("C23").Value = -3.834875614454537
$ws.Range("D23").Value = -4.855638079197258
$ws.Range("E23").Value = -5.032686349931778
$ws.Range("B24").Value = -0.001551833373864251
$ws.Range("C24").Value = -0.005756211124086755
$ws.Range("D24").Value = -0.006040061554452301
$ws.Range("E24").Value = -0.003238046131010818
$ws.Range("B25").Value = -1.091452332029469
$ws.Range("C25").Value = -2.016271703776692
$ws.Range("D25").Value = -2.014629965563369
$ws.Range("E25").Value = -1.037005651016436
